$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 134.1091918945312
    "B3" = 133.9417114257812
    "B4" = 135.1360931396484
    "B5" = 134.9424743652344
    "B6" = 132.00146484375
    "B7" = 131.8153076171875
    "B8" = 130.0806121826172
    "B9" = 129.9287872314453
    "B10" = 129.3352813720703
    "B11" = 129.2010192871094
    "B12" = 130.80517578125
    "B13" = 130.6780853271484
    "B14" = 140.7409362792969
    "B15" = 140.6166839599609
    "B16" = 159.1800689697266
    "B17" = 159.0573425292969
    "B18" = 172.2516479492188
    "B19" = 172.107666015625
    "B20" = 183.6344757080078
    "B21" = 183.4949798583984
    "B22" = 182.1724090576172
    "B23" = 182.0370025634766
    "B24" = 172.1752319335938
    "B25" = 172.0420989990234
    "B26" = 176.4307556152344
    "B27" = 176.2970428466797
    "B28" = 169.9020843505859
    "B29" = 169.7642517089844
    "B30" = 173.125732421875
    "B31" = 172.9806518554688
    "B32" = 187.027587890625
    "B33" = 186.873291015625
    "B34" = 215.2949523925781
    "B35" = 215.1307525634766
    "B36" = 239.4917907714844
    "B37" = 239.3180694580078
    "B38" = 213.1200408935547
    "B39" = 212.9378967285156
    "B40" = 176.8880310058594
    "B41" = 176.6986694335938
    "B42" = 156.5854949951172
    "B43" = 156.3897705078125
    "B44" = 140.8339538574219
    "B45" = 140.6324462890625
    "B46" = 126.9237670898438
    "B47" = 126.718017578125
    "B48" = 113.9836959838867
    "B49" = 113.7770156860352
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}
